$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Corrections to previously published rows (MV data update) ---
# Row 174 (Serie 01-05-2021)
$ws.Range("B174").Value = 908
$ws.Range("C174").Value = 796
$ws.Range("D174").Value = 3396
$ws.Range("I174").Value = 6107
$ws.Range("J174").Value = 5995
$ws.Range("N174").Value = 454
$ws.Range("O174").Value = 709
$ws.Range("P174").Value = 554
$ws.Range("Q174").Value = 1221

# Row 175 (Serie 01-06-2021)
$ws.Range("B175").Value = 580
$ws.Range("C175").Value = 468
$ws.Range("D175").Value = 3322
$ws.Range("I175").Value = 6288
$ws.Range("J175").Value = 6176
$ws.Range("N175").Value = 531
$ws.Range("O175").Value = 731
$ws.Range("P175").Value = 509
$ws.Range("Q175").Value = 1256

# Row 176 (Serie 01-07-2021)
$ws.Range("B176").Value = 245
$ws.Range("C176").Value = 133
$ws.Range("D176").Value = 2696
$ws.Range("I176").Value = 5371
$ws.Range("J176").Value = 5260
$ws.Range("N176").Value = 511
$ws.Range("O176").Value = 407
$ws.Range("P176").Value = 456
$ws.Range("Q176").Value = 1290

# --- New monthly row appended (Serie 01-08-2021) ---
# Enter the new period label as a text literal (via formula) and then
# convert it to a static value with PasteSpecial so Excel does not
# auto-convert the "01-08-2021" text into a date serial number.
$ws.Range("A177").Formula = "=""01-08-2021"""
$ws.Range("A177").Copy()
$ws.Range("A177").PasteSpecial(-4163)
$excel.CutCopyMode = 0

$ws.Range("B177").Value = 301
$ws.Range("C177").Value = 189
$ws.Range("D177").Value = 2770
$ws.Range("E177").Value = 2581
$ws.Range("F177").Value = 112
$ws.Range("G177").Value = 112
$ws.Range("H177").Value = 0
$ws.Range("I177").Value = 5463
$ws.Range("J177").Value = 5352
$ws.Range("K177").Value = 175
$ws.Range("L177").Value = 659
$ws.Range("M177").Value = 1706
$ws.Range("N177").Value = 671
$ws.Range("O177").Value = 388
$ws.Range("P177").Value = 430
$ws.Range("Q177").Value = 1323
$ws.Range("R177").Value = 112
